# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" column (D) for the row that is
# "Ready for handoff" (the fb0225ff-... source file, row 5) on both the
# zh-cn and de-de language sheets, recording the timestamp of the newly
# generated handoff.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-26 05:18:28"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-26 05:18:40"
